$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 1303
$ws.Range("F4").Value = 1105
$ws.Range("F5").Value = 990
$ws.Range("F7").Value = 546
$ws.Range("F12").Value = 277
$ws.Range("F14").Value = 85
$ws.Range("F15").Value = 660
$ws.Range("F16").Value = 150
$ws.Range("F20").Value = 326
$ws.Range("F21").Value = 128
$ws.Range("F22").Value = 656
$ws.Range("F23").Value = 24
$ws.Range("F24").Value = 637
$ws.Range("F27").Value = 864
$ws.Range("F28").Value = 307
$ws.Range("F29").Value = 144
$ws.Range("F30").Value = 35

# --- Sheet "演出" ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 11
$ws.Range("F7").Value = 247

# --- Sheet "本地生活" ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 304

# --- Sheet "全部类型" ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 304
$ws.Range("F4").Value = 1303
$ws.Range("F5").Value = 1105
$ws.Range("F6").Value = 990
$ws.Range("F8").Value = 546
$ws.Range("F14").Value = 277
$ws.Range("F16").Value = 85
$ws.Range("F17").Value = 660
$ws.Range("F18").Value = 150
$ws.Range("F24").Value = 11
$ws.Range("F25").Value = 326
$ws.Range("F27").Value = 247
$ws.Range("F28").Value = 247
$ws.Range("F29").Value = 128
$ws.Range("F30").Value = 656
$ws.Range("F31").Value = 24
$ws.Range("F32").Value = 637
$ws.Range("F35").Value = 864
$ws.Range("F36").Value = 307
$ws.Range("F39").Value = 144
$ws.Range("F40").Value = 35

$wb.Save()
